$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62, column A ("phone") was stored as text "51616172" and should become
# a genuine number 51616172 (same visible value, numeric type).
$ws.Range("A62").Value = 51616172

# New row 63: payment 51616176 (Cash), 2025-08-20T07:33:39
# Column A keeps the phone number as TEXT (matches how the id is stored
# elsewhere as an inline string), so force a text entry via the classic
# apostrophe-prefix trick, then strip the "quote prefix" style it leaves
# behind so no stray formatting is attached to the cell.
$ws.Range("A63").Formula = "'51616176"
$ws.Range("A63").Style = "Normal"

# B63 ("amount") and F63 ("discount_applied") are blank text cells, just
# like the other blank cells on this sheet (empty inline string, not a
# fully-cleared/empty cell). Use the same apostrophe trick with nothing
# after it to produce an empty text value, then strip the leftover style.
$ws.Range("B63").Formula = "'"
$ws.Range("B63").Style = "Normal"

$ws.Range("C63").Value = "Cash"
$ws.Range("D63").Value = "2025-08-20T07:33:39"
$ws.Range("E63").Value = 125

$ws.Range("F63").Formula = "'"
$ws.Range("F63").Style = "Normal"

$ws.Range("G63").Value = 125
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
